# data update 12 April 20
# Adds one new day (2020-04-12, serial 43933) of COVID-19 Bangladesh data
# to each of the three tracking sheets (Confirmed / Recoverd / Death),
# extending the running-total formula in column B down one row.

$wb = $excel.ActiveWorkbook

$wsConfirmed = $wb.Worksheets.Item("Confirmed")
$wsRecoverd  = $wb.Worksheets.Item("Recoverd")
$wsDeath     = $wb.Worksheets.Item("Death")

# ---------------------------------------------------------------------------
# Confirmed (sheet1): new row 37 -> 43933 | =SUM(B36+C37) | 139
# ---------------------------------------------------------------------------
$wsConfirmed.Range("A37").Value = 43933
$wsConfirmed.Range("A37").NumberFormat = $wsConfirmed.Range("A36").NumberFormat

$wsConfirmed.Range("B37").Formula = "=SUM(B36+C37)"
$wsConfirmed.Range("B37").HorizontalAlignment = $wsConfirmed.Range("B36").HorizontalAlignment
$wsConfirmed.Range("B37").VerticalAlignment = $wsConfirmed.Range("B36").VerticalAlignment

$wsConfirmed.Range("C37").Value = 139
$wsConfirmed.Range("C37").HorizontalAlignment = $wsConfirmed.Range("C36").HorizontalAlignment
$wsConfirmed.Range("C37").VerticalAlignment = $wsConfirmed.Range("C36").VerticalAlignment

# ---------------------------------------------------------------------------
# Recoverd (sheet2): new row 37 -> 43933 | =SUM(B36+C37) | 0
# ---------------------------------------------------------------------------
$wsRecoverd.Range("A37").Value = 43933
$wsRecoverd.Range("A37").NumberFormat = $wsRecoverd.Range("A36").NumberFormat

$wsRecoverd.Range("B37").Formula = "=SUM(B36+C37)"
$wsRecoverd.Range("B37").HorizontalAlignment = $wsRecoverd.Range("B36").HorizontalAlignment
$wsRecoverd.Range("B37").VerticalAlignment = $wsRecoverd.Range("B36").VerticalAlignment

$wsRecoverd.Range("C37").Value = 0
$wsRecoverd.Range("C37").HorizontalAlignment = $wsRecoverd.Range("C36").HorizontalAlignment
$wsRecoverd.Range("C37").VerticalAlignment = $wsRecoverd.Range("C36").VerticalAlignment

# ---------------------------------------------------------------------------
# Death (sheet3): new row 37 -> 43933 | =SUM(B36+C37) | 4
# ---------------------------------------------------------------------------
$wsDeath.Range("A37").Value = 43933
$wsDeath.Range("A37").NumberFormat = $wsDeath.Range("A36").NumberFormat

$wsDeath.Range("B37").Formula = "=SUM(B36+C37)"
$wsDeath.Range("B37").HorizontalAlignment = $wsDeath.Range("B36").HorizontalAlignment
$wsDeath.Range("B37").VerticalAlignment = $wsDeath.Range("B36").VerticalAlignment

$wsDeath.Range("C37").Value = 4
$wsDeath.Range("C37").HorizontalAlignment = $wsDeath.Range("C36").HorizontalAlignment
$wsDeath.Range("C37").VerticalAlignment = $wsDeath.Range("C36").VerticalAlignment

# ---------------------------------------------------------------------------
# Scroll/selection bookkeeping to mirror the author's on-screen state.
# Death is set last so it stays the active/selected tab (tabSelected),
# matching the workbook's saved activeTab.
# ---------------------------------------------------------------------------
$wsConfirmed.Range("J32").Select()
$excel.ActiveWindow.ScrollRow = 22
$excel.ActiveWindow.ScrollColumn = 1

$wsRecoverd.Range("D40").Select()
$excel.ActiveWindow.ScrollRow = 10
$excel.ActiveWindow.ScrollColumn = 1

$wsDeath.Range("D40").Select()
